# BigOp is PERFECT. Yay!
# Merge the two "SMALL" rows for tblTimeHistDetail into a single
# tblAdjustments/DeptNo row, rename the remaining "BIG" tblTimeHistDetail
# rows to tblTimeHistDetailOld, and highlight the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename tblTimeHistDetail -> tblTimeHistDetailOld for the BIG rows first (rows 6-12,
# before the row shift) so the new shared string is interned ahead of "tblAdjustments"
for ($r = 6; $r -le 12; $r++) {
    if ($ws.Cells.Item($r, 4).Value2 -eq "tblTimeHistDetail") {
        $ws.Cells.Item($r, 4).Value = "tblTimeHistDetailOld"
    }
}

# Row 3 becomes the merged/replacement row: SMALL / TimeCurrent / dbo / tblAdjustments / DeptNo / smallint / int
$ws.Range("A3").Value = "SMALL"
$ws.Range("B3").Value = "TimeCurrent"
$ws.Range("C3").Value = "dbo"
$ws.Range("D3").Value = "tblAdjustments"
$ws.Range("E3").Value = "DeptNo"
$ws.Range("F3").Value = "smallint"
$ws.Range("G3").Value = "int"

# Highlight B3:G3 with the new fill (Gold, Accent 4, Lighter 80% ~ theme accent4 / tint 0.8)
$ws.Range("B3:G3").Interior.Color = 13431551

# Shift rows 5-12 up to rows 4-11 (delete old row 4, the second SMALL row)
$ws.Rows("4").Delete()

# Column D width adjustment (widened to fit the longer "tblTimeHistDetailOld")
$ws.Columns("D").ColumnWidth = 17.15

# Selection moves to C17
$ws.Range("C17").Select()
